$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283, shifting existing rows 283-331 down to 284-332.
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new weekly price record.
$ws.Cells.Item(283, 1).Value = 8
$ws.Cells.Item(283, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(283, 3).Value = "Coquimbo"
$ws.Cells.Item(283, 4).Value = 45258
$ws.Cells.Item(283, 5).Value = 4
$ws.Cells.Item(283, 6).Value = 100112001
$ws.Cells.Item(283, 7).Value = "Berenjena"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 500
$ws.Cells.Item(283, 11).Value = 11000
$ws.Cells.Item(283, 12).Value = 12000
$ws.Cells.Item(283, 13).Value = 11500
$ws.Cells.Item(283, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(283, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(283, 16).Value = 230
$ws.Cells.Item(283, 17).Value = 50
$ws.Cells.Item(283, 18).Value = "Hortaliza"
